# Generate Report for Handoff
# Bumps the handoff timestamps for the six files that just got packaged for
# handoff, and marks them with the "ht" (handoff type) priority on the
# per-language localization-status sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows that correspond to the files handed off in this run.
$rows = @(7, 8, 11, 12, 13, 14)

foreach ($r in $rows) {
    # Overview!G<r> -> "Latest HO Xliff Generate Date"
    $overview.Cells.Item($r, 7).Value = "2016-08-28 10:23:54"

    # zh-cn!H<r> -> "Latest Handoff Datetime"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-28 10:23:50"

    # zh-cn!E<r> and de-de!E<r> -> "Priority" becomes "ht"
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 5).Value = "ht"
}
